# "The Latest results of the Arena"
# Applies the refreshed Arena leaderboard numbers, re-labels the rotated
# model names, (re)creates the results table, and updates the active
# selection/window state to match the latest authoring session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Window / view bookkeeping (best effort - some of this is read-only
#    chrome state in headless automation, but we still set it so that
#    a real Excel session would persist it).
# ---------------------------------------------------------------------
$excel.Height = 9180
$excel.ActiveWindow.Height = 9180

# ---------------------------------------------------------------------
# 2. The model lineup rotated: qwen/qwq-32b, open-r1/olympiccoder-32b and
#    google/gemma-3-27b-it swapped places in the standings.
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "google/gemma-3-27b-it:free"
$ws.Range("A5").Value = "qwen/qwq-32b:free"
$ws.Range("A6").Value = "open-r1/olympiccoder-32b:free"

# ---------------------------------------------------------------------
# 3. Refreshed ELO / timing / token / run statistics.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 2336
$ws.Range("C2").Value = 4.18
$ws.Range("D2").Value = 151
$ws.Range("E2").Value = 2

$ws.Range("B3").Value = 2058
$ws.Range("C3").Value = 4.41
$ws.Range("D3").Value = 66.5
$ws.Range("E3").Value = 2

$ws.Range("B4").Value = 1959
$ws.Range("C4").Value = 20.21
$ws.Range("D4").Value = 721.75
$ws.Range("E4").Value = 4

$ws.Range("B5").Value = 1836
$ws.Range("C5").Value = 19.76
$ws.Range("D5").Value = 666
$ws.Range("E5").Value = 1

$ws.Range("B6").Value = 1800
$ws.Range("C6").Value = 48.77
$ws.Range("D6").Value = 167.33
$ws.Range("E6").Value = 3

$ws.Range("B8").Value = 1687
$ws.Range("C8").Value = 25.14
$ws.Range("D8").Value = 204.5
$ws.Range("E8").Value = 2

$ws.Range("B9").Value = 1531
$ws.Range("C9").Value = 14.76
$ws.Range("D9").Value = 477.2
$ws.Range("E9").Value = 2

# ---------------------------------------------------------------------
# 4. Turn the results range into a proper table ("Tableau2"), sorted by
#    ELO descending, styled with TableStyleMedium9.
# ---------------------------------------------------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:F9"), [System.Type]::Missing, 1)
$tbl.Name = "Tableau2"
$tbl.TableStyle = "TableStyleMedium9"

$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($ws.Range("B2:B9"), [System.Type]::Missing, 2)
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# ---------------------------------------------------------------------
# 5. Final cursor position left by the author.
# ---------------------------------------------------------------------
$ws.Range("A6").Select()
